$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.75
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 4.75
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("Z3").Value = 17
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 19
$ws.Range("AK3").Value = 41
$ws.Range("AS3").Value = 251
$ws.Range("AZ3").Value = 81
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 3.6
$ws.Range("J5").Value = 3.4
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("Y5").Value = 10
$ws.Range("AC5").Value = 11
$ws.Range("AG5").Value = 201
$ws.Range("AH5").Value = 8.5
$ws.Range("AN5").Value = 4.75
$ws.Range("AO5").Value = 15
$ws.Range("AU5").Value = 7.5
$ws.Range("G7").Value = 2.67
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 2.6
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 2.02
$ws.Range("O7").Value = 1.24
$ws.Range("P7").Value = 3.3
$ws.Range("Q7").Value = 1.78
$ws.Range("R7").Value = 1.93
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.18
$ws.Range("X7").Value = 15.5
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 11
$ws.Range("AF7").Value = 40
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 9.75
$ws.Range("AI7").Value = 14.5
$ws.Range("AK7").Value = 32
$ws.Range("AM7").Value = 25
$ws.Range("AU7").Value = 6.2
$ws.Range("AV7").Value = 45
$ws.Range("AW7").Value = 4.65
$ws.Range("AY7").Value = 19
$ws.Range("BA7").Value = 80
$ws.Range("BC8").Value = 126
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 5.75
$ws.Range("K9").Value = 2.4
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("S9").Value = 1.3
$ws.Range("T9").Value = 3.4
$ws.Range("Z9").Value = 11
$ws.Range("AC9").Value = 13
$ws.Range("AD9").Value = 8
$ws.Range("AI9").Value = 34
$ws.Range("AO9").Value = 7.5
$ws.Range("AS9").Value = 101
$ws.Range("AT9").Value = 3.4
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 3.25
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 2.63
$ws.Range("L10").Value = 4.75
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("AC10").Value = 7.5
$ws.Range("AE10").Value = 19
$ws.Range("AK10").Value = 41
$ws.Range("AO10").Value = 11
$ws.Range("AQ10").Value = 41
$ws.Range("AX10").Value = 23
$ws.Range("BD10").Value = 151
$ws.Range("N11").Value = 8
$ws.Range("W12").Value = 6
$ws.Range("AN12").Value = 4
$ws.Range("AP12").Value = 29
$ws.Range("AQ12").Value = 51
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.98
